# Apply the FlashScore weekly-games update for 2024-10-28
# - Tweak several odds values in rows 3, 4 and 5
# - Replace row 6 (Peru Liga 1: Cusco vs Grau) with the USA MLS
#   (Seattle Sounders vs Houston Dynamo) match data, with refreshed odds
# - Remove the now-duplicate row 7, shrinking the sheet to A1:BD6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 odds tweaks ---
$ws.Range("G3").Value = 2.8
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 2.5
$ws.Range("L3").Value = 3.1
$ws.Range("Y3").Value = 11
$ws.Range("AA3").Value = 23
$ws.Range("AX3").Value = 13

# --- Row 4 odds tweaks ---
$ws.Range("G4").Value = 2.3
$ws.Range("I4").Value = 3.2
$ws.Range("J4").Value = 3.1
$ws.Range("K4").Value = 2.05
$ws.Range("L4").Value = 3.75
$ws.Range("U4").Value = 1.91
$ws.Range("V4").Value = 1.8
$ws.Range("W4").Value = 7
$ws.Range("AC4").Value = 8
$ws.Range("AM4").Value = 351
$ws.Range("BA4").Value = 81

# --- Row 5 odds tweaks ---
$ws.Range("Q5").Value = 2.15
$ws.Range("R5").Value = 1.67

# --- Row 6: replace match details with the USA - MLS game ---
$ws.Range("A6").Value = "CC5M2P9d"
$ws.Range("C6").Value = "21:50"
$ws.Range("D6").Value = "USA - MLS"
$ws.Range("E6").Value = "Seattle Sounders"
$ws.Range("F6").Value = "Houston Dynamo"

# --- Row 6: refreshed odds values ---
$ws.Range("G6").Value = 1.75
$ws.Range("H6").Value = 3.6
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 2.38
$ws.Range("L6").Value = 5.5
$ws.Range("O6").Value = 1.36
$ws.Range("P6").Value = 3
$ws.Range("Q6").Value = 2.15
$ws.Range("R6").Value = 1.67
$ws.Range("U6").Value = 2
$ws.Range("V6").Value = 1.75
$ws.Range("W6").Value = 6
$ws.Range("X6").Value = 7.5
$ws.Range("Z6").Value = 13
$ws.Range("AA6").Value = 15
$ws.Range("AB6").Value = 34
$ws.Range("AD6").Value = 6.5
$ws.Range("AE6").Value = 19
$ws.Range("AF6").Value = 67
$ws.Range("AH6").Value = 23
$ws.Range("AJ6").Value = 51
$ws.Range("AK6").Value = 41
$ws.Range("AM6").Value = 451
$ws.Range("AN6").Value = 3.6
$ws.Range("AO6").Value = 9.5
$ws.Range("AQ6").Value = 34
$ws.Range("AS6").Value = 201
$ws.Range("AU6").Value = 9
$ws.Range("AV6").Value = 67
$ws.Range("AW6").Value = 6.5
$ws.Range("AX6").Value = 29
$ws.Range("AY6").Value = 34
$ws.Range("AZ6").Value = 101
$ws.Range("BA6").Value = 126
$ws.Range("BB6").Value = 301

# --- Remove the old duplicate row 7 (its data now lives in row 6) ---
$ws.Rows.Item(7).Delete()
